# Delete the post at row 142 ("「失敗は成功の始まり」...").
# Deleting the entire row shifts every subsequent row up by one,
# which also updates the used range from A1:C290 to A1:C289.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Rows.Item(142).Delete()
